$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 162.399297
$ws.Cells.Item(2, 8).Value = 487.197891
$ws.Cells.Item(2, 9).Value = 0.3910371682630009
$ws.Cells.Item(2, 10).Value = 0.3910371682630009
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 23.18520366666667
$ws.Cells.Item(2, 14).Value = 69.555611
$ws.Cells.Item(2, 15).Value = 0.4216200689608106
$ws.Cells.Item(2, 16).Value = 0.4216200689608105
$ws.Cells.Item(2, 17).Value = 3765.260776268489
$ws.Cells.Item(2, 18).Value = 33887.3469864164
$ws.Cells.Item(2, 19).Value = 0.1648691178492865
$ws.Cells.Item(2, 20).Value = 0.1648691178492865

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 162.399297
$ws.Cells.Item(3, 8).Value = 487.197891
$ws.Cells.Item(3, 9).Value = 0.3910371682630009
$ws.Cells.Item(3, 10).Value = 0.3910371682630009
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 11.56543033333333
$ws.Cells.Item(3, 14).Value = 34.696291
$ws.Cells.Item(3, 15).Value = 0.210315924104302
$ws.Cells.Item(3, 16).Value = 0.2103159241043019
$ws.Cells.Item(3, 17).Value = 1878.217755635809
$ws.Cells.Item(3, 18).Value = 16903.95980072228
$ws.Cells.Item(3, 19).Value = 0.08224134340236247
$ws.Cells.Item(3, 20).Value = 0.08224134340236246

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 162.399297
$ws.Cells.Item(4, 8).Value = 487.197891
$ws.Cells.Item(4, 9).Value = 0.3910371682630009
$ws.Cells.Item(4, 10).Value = 0.3910371682630009
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 2.096289333333333
$ws.Cells.Item(4, 14).Value = 6.288868
$ws.Cells.Item(4, 15).Value = 0.03812076296541245
$ws.Cells.Item(4, 16).Value = 0.03812076296541244
$ws.Cells.Item(4, 17).Value = 340.435914041932
$ws.Cells.Item(4, 18).Value = 3063.923226377388
$ws.Cells.Item(4, 19).Value = 0.01490663520201996
$ws.Cells.Item(4, 20).Value = 0.01490663520201996

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 162.399297
$ws.Cells.Item(5, 8).Value = 487.197891
$ws.Cells.Item(5, 9).Value = 0.3910371682630009
$ws.Cells.Item(5, 10).Value = 0.3910371682630009
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 18.14382633333333
$ws.Cells.Item(5, 14).Value = 54.431479
$ws.Cells.Item(5, 15).Value = 0.3299432439694752
$ws.Cells.Item(5, 16).Value = 0.3299432439694752
$ws.Cells.Item(5, 17).Value = 2946.544641423421
$ws.Cells.Item(5, 18).Value = 26518.90177281079
$ws.Cells.Item(5, 19).Value = 0.129020071809332
$ws.Cells.Item(5, 20).Value = 0.129020071809332

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 65.41736466666667
$ws.Cells.Item(6, 8).Value = 196.252094
$ws.Cells.Item(6, 9).Value = 0.1575168212364948
$ws.Cells.Item(6, 10).Value = 0.1575168212364948
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 23.18520366666667
$ws.Cells.Item(6, 14).Value = 69.555611
$ws.Cells.Item(6, 15).Value = 0.4216200689608106
$ws.Cells.Item(6, 16).Value = 0.4216200689608105
$ws.Cells.Item(6, 17).Value = 1516.714923133271
$ws.Cells.Item(6, 18).Value = 13650.43430819943
$ws.Cells.Item(6, 19).Value = 0.06641225303221862
$ws.Cells.Item(6, 20).Value = 0.06641225303221861

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 65.41736466666667
$ws.Cells.Item(7, 8).Value = 196.252094
$ws.Cells.Item(7, 9).Value = 0.1575168212364948
$ws.Cells.Item(7, 10).Value = 0.1575168212364948
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 11.56543033333333
$ws.Cells.Item(7, 14).Value = 34.696291
$ws.Cells.Item(7, 15).Value = 0.210315924104302
$ws.Cells.Item(7, 16).Value = 0.2103159241043019
$ws.Cells.Item(7, 17).Value = 756.579973642595
$ws.Cells.Item(7, 18).Value = 6809.219762783355
$ws.Cells.Item(7, 19).Value = 0.03312829582032555
$ws.Cells.Item(7, 20).Value = 0.03312829582032555

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 65.41736466666667
$ws.Cells.Item(8, 8).Value = 196.252094
$ws.Cells.Item(8, 9).Value = 0.1575168212364948
$ws.Cells.Item(8, 10).Value = 0.1575168212364948
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 2.096289333333333
$ws.Cells.Item(8, 14).Value = 6.288868
$ws.Cells.Item(8, 15).Value = 0.03812076296541245
$ws.Cells.Item(8, 16).Value = 0.03812076296541244
$ws.Cells.Item(8, 17).Value = 137.1337237655102
$ws.Cells.Item(8, 18).Value = 1234.203513889592
$ws.Cells.Item(8, 19).Value = 0.006004661405421666
$ws.Cells.Item(8, 20).Value = 0.006004661405421664

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 65.41736466666667
$ws.Cells.Item(9, 8).Value = 196.252094
$ws.Cells.Item(9, 9).Value = 0.1575168212364948
$ws.Cells.Item(9, 10).Value = 0.1575168212364948
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 18.14382633333333
$ws.Cells.Item(9, 14).Value = 54.431479
$ws.Cells.Item(9, 15).Value = 0.3299432439694752
$ws.Cells.Item(9, 16).Value = 0.3299432439694752
$ws.Cells.Item(9, 17).Value = 1186.921303696336
$ws.Cells.Item(9, 18).Value = 10682.29173326703
$ws.Cells.Item(9, 19).Value = 0.05197161097852902
$ws.Cells.Item(9, 20).Value = 0.05197161097852902

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 126.3069433333333
$ws.Cells.Item(10, 8).Value = 378.92083
$ws.Cells.Item(10, 9).Value = 0.3041313008456065
$ws.Cells.Item(10, 10).Value = 0.3041313008456065
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 23.18520366666667
$ws.Cells.Item(10, 14).Value = 69.555611
$ws.Cells.Item(10, 15).Value = 0.4216200689608106
$ws.Cells.Item(10, 16).Value = 0.4216200689608105
$ws.Cells.Item(10, 17).Value = 2928.452205697459
$ws.Cells.Item(10, 18).Value = 26356.06985127713
$ws.Cells.Item(10, 19).Value = 0.1282278600356656
$ws.Cells.Item(10, 20).Value = 0.1282278600356656

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 126.3069433333333
$ws.Cells.Item(11, 8).Value = 378.92083
$ws.Cells.Item(11, 9).Value = 0.3041313008456065
$ws.Cells.Item(11, 10).Value = 0.3041313008456065
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 11.56543033333333
$ws.Cells.Item(11, 14).Value = 34.696291
$ws.Cells.Item(11, 15).Value = 0.210315924104302
$ws.Cells.Item(11, 16).Value = 0.2103159241043019
$ws.Cells.Item(11, 17).Value = 1460.794153737948
$ws.Cells.Item(11, 18).Value = 13147.14738364153
$ws.Cells.Item(11, 19).Value = 0.0639636555863872
$ws.Cells.Item(11, 20).Value = 0.0639636555863872

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 126.3069433333333
$ws.Cells.Item(12, 8).Value = 378.92083
$ws.Cells.Item(12, 9).Value = 0.3041313008456065
$ws.Cells.Item(12, 10).Value = 0.3041313008456065
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 2.096289333333333
$ws.Cells.Item(12, 14).Value = 6.288868
$ws.Cells.Item(12, 15).Value = 0.03812076296541245
$ws.Cells.Item(12, 16).Value = 0.03812076296541244
$ws.Cells.Item(12, 17).Value = 264.7758980356045
$ws.Cells.Item(12, 18).Value = 2382.98308232044
$ws.Cells.Item(12, 19).Value = 0.01159371722989791
$ws.Cells.Item(12, 20).Value = 0.01159371722989791

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 126.3069433333333
$ws.Cells.Item(13, 8).Value = 378.92083
$ws.Cells.Item(13, 9).Value = 0.3041313008456065
$ws.Cells.Item(13, 10).Value = 0.3041313008456065
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 18.14382633333333
$ws.Cells.Item(13, 14).Value = 54.431479
$ws.Cells.Item(13, 15).Value = 0.3299432439694752
$ws.Cells.Item(13, 16).Value = 0.3299432439694752
$ws.Cells.Item(13, 17).Value = 2291.691244534175
$ws.Cells.Item(13, 18).Value = 20625.22120080757
$ws.Cells.Item(13, 19).Value = 0.1003460679936558
$ws.Cells.Item(13, 20).Value = 0.1003460679936558

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 61.180387
$ws.Cells.Item(14, 8).Value = 183.541161
$ws.Cells.Item(14, 9).Value = 0.1473147096548978
$ws.Cells.Item(14, 10).Value = 0.1473147096548978
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 23.18520366666667
$ws.Cells.Item(14, 14).Value = 69.555611
$ws.Cells.Item(14, 15).Value = 0.4216200689608106
$ws.Cells.Item(14, 16).Value = 0.4216200689608105
$ws.Cells.Item(14, 17).Value = 1418.479733000486
$ws.Cells.Item(14, 18).Value = 12766.31759700437
$ws.Cells.Item(14, 19).Value = 0.06211083804363979
$ws.Cells.Item(14, 20).Value = 0.06211083804363979

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 61.180387
$ws.Cells.Item(15, 8).Value = 183.541161
$ws.Cells.Item(15, 9).Value = 0.1473147096548978
$ws.Cells.Item(15, 10).Value = 0.1473147096548978
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 11.56543033333333
$ws.Cells.Item(15, 14).Value = 34.696291
$ws.Cells.Item(15, 15).Value = 0.210315924104302
$ws.Cells.Item(15, 16).Value = 0.2103159241043019
$ws.Cells.Item(15, 17).Value = 707.5775036148723
$ws.Cells.Item(15, 18).Value = 6368.197532533851
$ws.Cells.Item(15, 19).Value = 0.03098262929522676
$ws.Cells.Item(15, 20).Value = 0.03098262929522677

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 61.180387
$ws.Cells.Item(16, 8).Value = 183.541161
$ws.Cells.Item(16, 9).Value = 0.1473147096548978
$ws.Cells.Item(16, 10).Value = 0.1473147096548978
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 2.096289333333333
$ws.Cells.Item(16, 14).Value = 6.288868
$ws.Cells.Item(16, 15).Value = 0.03812076296541245
$ws.Cells.Item(16, 16).Value = 0.03812076296541244
$ws.Cells.Item(16, 17).Value = 128.2517926773053
$ws.Cells.Item(16, 18).Value = 1154.266134095748
$ws.Cells.Item(16, 19).Value = 0.005615749128072916
$ws.Cells.Item(16, 20).Value = 0.005615749128072915

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 61.180387
$ws.Cells.Item(17, 8).Value = 183.541161
$ws.Cells.Item(17, 9).Value = 0.1473147096548978
$ws.Cells.Item(17, 10).Value = 0.1473147096548978
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 18.14382633333333
$ws.Cells.Item(17, 14).Value = 54.431479
$ws.Cells.Item(17, 15).Value = 0.3299432439694752
$ws.Cells.Item(17, 16).Value = 0.3299432439694752
$ws.Cells.Item(17, 17).Value = 1110.046316734124
$ws.Cells.Item(17, 18).Value = 9990.416850607118
$ws.Cells.Item(17, 19).Value = 0.04860549318795834
$ws.Cells.Item(17, 20).Value = 0.04860549318795835
